$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 49 - this shifts the existing rows 49:175 down to 50:176
# (matches the Excel "insert shifting cells down" behaviour, carrying the row formatting
# from the surrounding rows, same as the author's manual row insert before typing new data).
$ws.Rows(49).Insert()

# Fill in the new weekly data record for row 49.
$ws.Range("A49").Value = 5
$ws.Range("B49").Value = "Macroferia Regional de Talca"
$ws.Range("C49").Value = "Maule"
$ws.Range("D49").Value = 44953
$ws.Range("E49").Value = 7
$ws.Range("F49").Value = 100112030
$ws.Range("G49").Value = "Poroto granado"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 100
$ws.Range("K49").Value = 40000
$ws.Range("L49").Value = 40000
$ws.Range("M49").Value = 40000
$ws.Range("N49").Value = "`$/saco 25 kilos"
$ws.Range("O49").Value = "Región del Maule"
$ws.Range("P49").Value = 1600
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"
